$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update accuracy values in column B (rows changed per diff)
$ws.Range("B2").Value = 0.375
$ws.Range("B3").Value = 0.203125
$ws.Range("B4").Value = 0.171875
$ws.Range("B5").Value = 0.15625
$ws.Range("B6").Value = 0.15625
$ws.Range("B7").Value = 0.203125
$ws.Range("B8").Value = 0.15625
$ws.Range("B9").Value = 0.15625
$ws.Range("B10").Value = 0.125
$ws.Range("B11").Value = 0.140625
$ws.Range("B12").Value = 0.125
$ws.Range("B13").Value = 0.15625
$ws.Range("B14").Value = 0.140625
$ws.Range("B15").Value = 0.109375
$ws.Range("B16").Value = 0.140625
$ws.Range("B17").Value = 0.203125
$ws.Range("B20").Value = 0.125
$ws.Range("B21").Value = 0.28125
$ws.Range("B22").Value = 0.140625
$ws.Range("B23").Value = 0.15625
$ws.Range("B27").Value = 0.125
$ws.Range("B28").Value = 0.125
$ws.Range("B29").Value = 0.125
$ws.Range("B30").Value = 0.125
$ws.Range("B31").Value = 0.125
$ws.Range("B36").Value = 0.140625
$ws.Range("B43").Value = 0.125
$ws.Range("B44").Value = 0.125
$ws.Range("B45").Value = 0.125
$ws.Range("B46").Value = 0.125
$ws.Range("B47").Value = 0.125
$ws.Range("B48").Value = 0.125
$ws.Range("B49").Value = 0.125
$ws.Range("B50").Value = 0.125
$ws.Range("B51").Value = 0.125
$ws.Range("B52").Value = 0.125
$ws.Range("B53").Value = 0.125
$ws.Range("B54").Value = 0.125
$ws.Range("B55").Value = 0.125
$ws.Range("B56").Value = 0.125
$ws.Range("B106").Value = 0.09375
$ws.Range("B108").Value = 0.078125
$ws.Range("B109").Value = 0.046875
$ws.Range("B111").Value = 0.0625
$ws.Range("B112").Value = 0.09375
$ws.Range("B113").Value = 0.046875
$ws.Range("B114").Value = 0.09375
$ws.Range("B115").Value = 0.03125
$ws.Range("B116").Value = 0
$ws.Range("B117").Value = 0.078125
$ws.Range("B118").Value = 0.03278688524590164
# Update DisplayOutputs object address string in column A (rows 102-118)
$ws.Range("A102").Value = "<__main__.DisplayOutputs object at 0x7f6040058b20>"
$ws.Range("A103").Value = "<__main__.DisplayOutputs object at 0x7f6040058b20>"
$ws.Range("A104").Value = "<__main__.DisplayOutputs object at 0x7f6040058b20>"
$ws.Range("A105").Value = "<__main__.DisplayOutputs object at 0x7f6040058b20>"
$ws.Range("A106").Value = "<__main__.DisplayOutputs object at 0x7f6040058b20>"
$ws.Range("A107").Value = "<__main__.DisplayOutputs object at 0x7f6040058b20>"
$ws.Range("A108").Value = "<__main__.DisplayOutputs object at 0x7f6040058b20>"
$ws.Range("A109").Value = "<__main__.DisplayOutputs object at 0x7f6040058b20>"
$ws.Range("A110").Value = "<__main__.DisplayOutputs object at 0x7f6040058b20>"
$ws.Range("A111").Value = "<__main__.DisplayOutputs object at 0x7f6040058b20>"
$ws.Range("A112").Value = "<__main__.DisplayOutputs object at 0x7f6040058b20>"
$ws.Range("A113").Value = "<__main__.DisplayOutputs object at 0x7f6040058b20>"
$ws.Range("A114").Value = "<__main__.DisplayOutputs object at 0x7f6040058b20>"
$ws.Range("A115").Value = "<__main__.DisplayOutputs object at 0x7f6040058b20>"
$ws.Range("A116").Value = "<__main__.DisplayOutputs object at 0x7f6040058b20>"
$ws.Range("A117").Value = "<__main__.DisplayOutputs object at 0x7f6040058b20>"
$ws.Range("A118").Value = "<__main__.DisplayOutputs object at 0x7f6040058b20>"